$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Rtn4"
$ws.Range("C2").Value = "Cntnap1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 39.66867433333334
$ws.Range("H2").Value = 119.006023
$ws.Range("I2").Value = 0.154574216411057
$ws.Range("J2").Value = 0.1545742164110569
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 1.067233
$ws.Range("N2").Value = 3.201699
$ws.Range("O2").Value = 0.23609279392832
$ws.Range("P2").Value = 0.23609279392832
$ws.Range("Q2").Value = 42.33571831478634
$ws.Range("R2").Value = 381.021464833077
$ws.Range("S2").Value = 0.03649385862176721
$ws.Range("T2").Value = 0.0364938586217672

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Rtn4"
$ws.Range("C3").Value = "Cntnap1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 39.66867433333334
$ws.Range("H3").Value = 119.006023
$ws.Range("I3").Value = 0.154574216411057
$ws.Range("J3").Value = 0.1545742164110569
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 2.034898666666666
$ws.Range("N3").Value = 6.104696
$ws.Range("O3").Value = 0.4501593481220562
$ws.Range("P3").Value = 0.4501593481220562
$ws.Range("Q3").Value = 80.72173250933422
$ws.Range("R3").Value = 726.4955925840079
$ws.Range("S3").Value = 0.06958302849607904
$ws.Range("T3").Value = 0.06958302849607903

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Rtn4"
$ws.Range("C4").Value = "Cntnap1"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 39.66867433333334
$ws.Range("H4").Value = 119.006023
$ws.Range("I4").Value = 0.154574216411057
$ws.Range("J4").Value = 0.1545742164110569
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.2589606666666667
$ws.Range("N4").Value = 0.776882
$ws.Range("O4").Value = 0.0572871597025895
$ws.Range("P4").Value = 0.0572871597025895
$ws.Range("Q4").Value = 10.27262635114289
$ws.Range("R4").Value = 92.45363716028599
$ws.Range("S4").Value = 0.00885511782144285
$ws.Range("T4").Value = 0.008855117821442849

$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Rtn4"
$ws.Range("C5").Value = "Cntnap1"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 39.66867433333334
$ws.Range("H5").Value = 119.006023
$ws.Range("I5").Value = 0.154574216411057
$ws.Range("J5").Value = 0.1545742164110569
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.159304
$ws.Range("N5").Value = 3.477912
$ws.Range("O5").Value = 0.2564606982470342
$ws.Range("P5").Value = 0.2564606982470342
$ws.Range("Q5").Value = 45.98805282933066
$ws.Range("R5").Value = 413.892475463976
$ws.Range("S5").Value = 0.03964221147176784
$ws.Range("T5").Value = 0.03964221147176784

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Rtn4"
$ws.Range("C6").Value = "Cntnap1"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 57.66057933333332
$ws.Range("H6").Value = 172.981738
$ws.Range("I6").Value = 0.2246820449144221
$ws.Range("J6").Value = 0.2246820449144221
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 1.067233
$ws.Range("N6").Value = 3.201699
$ws.Range("O6").Value = 0.23609279392832
$ws.Range("P6").Value = 0.23609279392832
$ws.Range("Q6").Value = 61.53727306365133
$ws.Range("R6").Value = 553.835457572862
$ws.Range("S6").Value = 0.0530458117293742
$ws.Range("T6").Value = 0.05304581172937418

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Rtn4"
$ws.Range("C7").Value = "Cntnap1"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 57.66057933333332
$ws.Range("H7").Value = 172.981738
$ws.Range("I7").Value = 0.2246820449144221
$ws.Range("J7").Value = 0.2246820449144221
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 2.034898666666666
$ws.Range("N7").Value = 6.104696
$ws.Range("O7").Value = 0.4501593481220562
$ws.Range("P7").Value = 0.4501593481220562
$ws.Range("Q7").Value = 117.3334360046275
$ws.Range("R7").Value = 1056.000924041648
$ws.Range("S7").Value = 0.1011427228734068
$ws.Range("T7").Value = 0.1011427228734068

$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Rtn4"
$ws.Range("C8").Value = "Cntnap1"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 57.66057933333332
$ws.Range("H8").Value = 172.981738
$ws.Range("I8").Value = 0.2246820449144221
$ws.Range("J8").Value = 0.2246820449144221
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.2589606666666667
$ws.Range("N8").Value = 0.776882
$ws.Range("O8").Value = 0.0572871597025895
$ws.Range("P8").Value = 0.0572871597025895
$ws.Range("Q8").Value = 14.93182206454622
$ws.Range("R8").Value = 134.386398580916
$ws.Range("S8").Value = 0.01287139618931689
$ws.Range("T8").Value = 0.01287139618931688

$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Rtn4"
$ws.Range("C9").Value = "Cntnap1"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 57.66057933333332
$ws.Range("H9").Value = 172.981738
$ws.Range("I9").Value = 0.2246820449144221
$ws.Range("J9").Value = 0.2246820449144221
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.159304
$ws.Range("N9").Value = 3.477912
$ws.Range("O9").Value = 0.2564606982470342
$ws.Range("P9").Value = 0.2564606982470342
$ws.Range("Q9").Value = 66.84614026345065
$ws.Range("R9").Value = 601.6152623710559
$ws.Range("S9").Value = 0.05762211412232419
$ws.Range("T9").Value = 0.05762211412232419

$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "Rtn4"
$ws.Range("C10").Value = "Cntnap1"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 64.993678
$ws.Range("H10").Value = 194.981034
$ws.Range("I10").Value = 0.2532564301015895
$ws.Range("J10").Value = 0.2532564301015895
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 1.067233
$ws.Range("N10").Value = 3.201699
$ws.Range("O10").Value = 0.23609279392832
$ws.Range("P10").Value = 0.23609279392832
$ws.Range("Q10").Value = 69.36339795297401
$ws.Range("R10").Value = 624.2705815767661
$ws.Range("S10").Value = 0.05979201816299656
$ws.Range("T10").Value = 0.05979201816299655

$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "Rtn4"
$ws.Range("C11").Value = "Cntnap1"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 64.993678
$ws.Range("H11").Value = 194.981034
$ws.Range("I11").Value = 0.2532564301015895
$ws.Range("J11").Value = 0.2532564301015895
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 2.034898666666666
$ws.Range("N11").Value = 6.104696
$ws.Range("O11").Value = 0.4501593481220562
$ws.Range("P11").Value = 0.4501593481220562
$ws.Range("Q11").Value = 132.2555487039627
$ws.Range("R11").Value = 1190.299938335664
$ws.Range("S11").Value = 0.1140057494822506
$ws.Range("T11").Value = 0.1140057494822506

$ws.Range("A12").Value = "M2"
$ws.Range("B12").Value = "Rtn4"
$ws.Range("C12").Value = "Cntnap1"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 64.993678
$ws.Range("H12").Value = 194.981034
$ws.Range("I12").Value = 0.2532564301015895
$ws.Range("J12").Value = 0.2532564301015895
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 0.2589606666666667
$ws.Range("N12").Value = 0.776882
$ws.Range("O12").Value = 0.0572871597025895
$ws.Range("P12").Value = 0.0572871597025895
$ws.Range("Q12").Value = 16.83080618399867
$ws.Range("R12").Value = 151.477255655988
$ws.Range("S12").Value = 0.01450834155693745
$ws.Range("T12").Value = 0.01450834155693745

$ws.Range("A13").Value = "M2"
$ws.Range("B13").Value = "Rtn4"
$ws.Range("C13").Value = "Cntnap1"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 64.993678
$ws.Range("H13").Value = 194.981034
$ws.Range("I13").Value = 0.2532564301015895
$ws.Range("J13").Value = 0.2532564301015895
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 1.159304
$ws.Range("N13").Value = 3.477912
$ws.Range("O13").Value = 0.2564606982470342
$ws.Range("P13").Value = 0.2564606982470342
$ws.Range("Q13").Value = 75.347430880112
$ws.Range("R13").Value = 678.126877921008
$ws.Range("S13").Value = 0.06495032089940485
$ws.Range("T13").Value = 0.06495032089940485

$ws.Range("A14").Value = "sCs"
$ws.Range("B14").Value = "Rtn4"
$ws.Range("C14").Value = "Cntnap1"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 94.308965
$ws.Range("H14").Value = 282.926895
$ws.Range("I14").Value = 0.3674873085729315
$ws.Range("J14").Value = 0.3674873085729314
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 0.6666666666666666
$ws.Range("M14").Value = 1.067233
$ws.Range("N14").Value = 3.201699
$ws.Range("O14").Value = 0.23609279392832
$ws.Range("P14").Value = 0.23609279392832
$ws.Range("Q14").Value = 100.649639643845
$ws.Range("R14").Value = 905.846756794605
$ws.Range("S14").Value = 0.08676110541418207
$ws.Range("T14").Value = 0.08676110541418204

$ws.Range("A15").Value = "sCs"
$ws.Range("B15").Value = "Rtn4"
$ws.Range("C15").Value = "Cntnap1"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 94.308965
$ws.Range("H15").Value = 282.926895
$ws.Range("I15").Value = 0.3674873085729315
$ws.Range("J15").Value = 0.3674873085729314
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 2.034898666666666
$ws.Range("N15").Value = 6.104696
$ws.Range("O15").Value = 0.4501593481220562
$ws.Range("P15").Value = 0.4501593481220562
$ws.Range("Q15").Value = 191.9091871332133
$ws.Range("R15").Value = 1727.18268419892
$ws.Range("S15").Value = 0.1654278472703198
$ws.Range("T15").Value = 0.1654278472703197

$ws.Range("A16").Value = "sCs"
$ws.Range("B16").Value = "Rtn4"
$ws.Range("C16").Value = "Cntnap1"
$ws.Range("D16").Value = "M2"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 94.308965
$ws.Range("H16").Value = 282.926895
$ws.Range("I16").Value = 0.3674873085729315
$ws.Range("J16").Value = 0.3674873085729314
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 0.2589606666666667
$ws.Range("N16").Value = 0.776882
$ws.Range("O16").Value = 0.0572871597025895
$ws.Range("P16").Value = 0.0572871597025895
$ws.Range("Q16").Value = 24.42231244904334
$ws.Range("R16").Value = 219.80081204139
$ws.Range("S16").Value = 0.02105230413489232
$ws.Range("T16").Value = 0.02105230413489231

$ws.Range("A17").Value = "sCs"
$ws.Range("B17").Value = "Rtn4"
$ws.Range("C17").Value = "Cntnap1"
$ws.Range("D17").Value = "sCs"
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 94.308965
$ws.Range("H17").Value = 282.926895
$ws.Range("I17").Value = 0.3674873085729315
$ws.Range("J17").Value = 0.3674873085729314
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 1.159304
$ws.Range("N17").Value = 3.477912
$ws.Range("O17").Value = 0.2564606982470342
$ws.Range("P17").Value = 0.2564606982470342
$ws.Range("Q17").Value = 109.33276036036
$ws.Range("R17").Value = 983.99484324324
$ws.Range("S17").Value = 0.09424605175353733
$ws.Range("T17").Value = 0.09424605175353731
